$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (can look numeric, e.g. "579.53" or "69.486.31");
# force text format first so Excel does not coerce these into numeric cells
# and lose formatting (trailing zeros, multi-dot thousands separators, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.486.31"
$ws.Range("E2").Value = "  -3.91%  "
$ws.Range("D3").Value = "2.506.51"
$ws.Range("E3").Value = "  -5.39%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "579.53"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "167.17"
$ws.Range("E6").Value = "  -4.52%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "2.504.85"
$ws.Range("E9").Value = "  -5.40%  "
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -6.92%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").Value = "2.964.48"
$ws.Range("E14").Value = "  -5.40%  "
$ws.Range("D15").Value = "69.402.58"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("E16").Value = "  -5.57%  "
$ws.Range("D17").Value = "24.94"
$ws.Range("E17").Value = "  -4.20%  "
$ws.Range("D18").Value = "2.508.29"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  -6.92%  "
$ws.Range("D20").Value = "7.81"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "351.17"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").Value = "3.96"
$ws.Range("E22").Value = "  -5.02%  "
$ws.Range("D23").Value = "1.97"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "69.19"
$ws.Range("E25").Value = "  -3.31%  "
$ws.Range("E26").Value = "  -5.58%  "
$ws.Range("E27").Value = "  -6.59%  "
$ws.Range("D28").Value = "2.634.93"
$ws.Range("E28").Value = "  -5.43%  "
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "0.0₃0905"
$ws.Range("E30").Value = "  -5.28%  "
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").Value = "479.67"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").Value = "1.29"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "153.19"
$ws.Range("E37").Value = "  -5.38%  "
$ws.Range("D38").Value = "18.88"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "18.58"
$ws.Range("E39").Value = "  -4.12%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  -6.11%  "
$ws.Range("D44").Value = "1.18"
$ws.Range("E44").Value = "  -13.06%  "
$ws.Range("E45").Value = "  -8.55%  "
$ws.Range("D46").Value = "38.16"
$ws.Range("D47").Value = "144.01"
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("D48").Value = "3.55"
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D49").Value = "0.530"
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("E50").Value = "  -5.17%  "
$ws.Range("D51").Value = "0.0731"
$ws.Range("E51").Value = "  -2.31%  "
